# Update countries & provincias Spain
# Refresh the "Datos actualizados a ..." timestamp string (always the last
# data row, row 1, column A) and apply the latest COVID country-stats
# snapshot: several rows changed their country label because the
# underlying ranking re-sorted (Peru overtook Colombia, Jamaica overtook
# Birmania, Lesoto overtook Liberia, Montserrat overtook Islas Malvinas),
# plus a batch of updated case counts for this and other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 02:54"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("B4").Value = 6925589
$ws.Range("C4").Value = 50993
$ws.Range("D4").Value = 4191752
$ws.Range("E4").Value = 2530697
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 927
$ws.Range("H4").Value = 203140

# --- Row 8: now Peru (was Colombia) ------------------------------------
$ws.Range("A8").Value = "Peru"
$ws.Range("B8").Value = 756412
$ws.Range("C8").Value = 6314
$ws.Range("D8").Value = 600795
$ws.Range("E8").Value = 124334
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 137
$ws.Range("H8").Value = 31283

# --- Row 9: now Colombia (was Peru) ------------------------------------
$ws.Range("A9").Value = "Colombia"
$ws.Range("B9").Value = 750471
$ws.Range("C9").Value = 6526
$ws.Range("D9").Value = 621521
$ws.Range("E9").Value = 105100
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 185
$ws.Range("H9").Value = 23850

# --- Row 13: Argentina ---------------------------------------------------
$ws.Range("B13").Value = 613658
$ws.Range("C13").Value = 11945
$ws.Range("E13").Value = 133716
$ws.Range("G13").Value = 196
$ws.Range("H13").Value = 12656

# --- Row 36: Panama --------------------------------------------------
$ws.Range("B36").Value = 104879
$ws.Range("C36").Value = 741
$ws.Range("D36").Value = 79093
$ws.Range("E36").Value = 23557
$ws.Range("G36").Value = 16
$ws.Range("H36").Value = 2229

# --- Row 74: Paraguay --------------------------------------------------
$ws.Range("B74").Value = 32127
$ws.Range("C74").Value = 1014
$ws.Range("D74").Value = 16921
$ws.Range("E74").Value = 14595
$ws.Range("G74").Value = 27
$ws.Range("H74").Value = 611

# --- Row 125: now Jamaica (was Birmania) --------------------------------
$ws.Range("A125").Value = "Jamaica"
$ws.Range("B125").Value = 4571
$ws.Range("C125").Value = 197
$ws.Range("D125").Value = 1264
$ws.Range("E125").Value = 3252
$ws.Range("G125").Value = 4
$ws.Range("H125").Value = 55

# --- Row 126: now Birmania (was Jamaica) --------------------------------
$ws.Range("A126").Value = "Birmania"
$ws.Range("B126").Value = 4467
$ws.Range("C126").Value = 424
$ws.Range("D126").Value = 1130
$ws.Range("E126").Value = 3267
$ws.Range("G126").Value = 10
$ws.Range("H126").Value = 70

# --- Row 162: now Lesoto (was Liberia) ----------------------------------
$ws.Range("A162").Value = "Lesoto"
$ws.Range("B162").Value = 1390
$ws.Range("C162").Value = 63
$ws.Range("D162").Value = 754
$ws.Range("E162").Value = 603
$ws.Range("H162").Value = 33

# --- Row 163: now Liberia (was Lesoto) ----------------------------------
$ws.Range("A163").Value = "Liberia"
$ws.Range("B163").Value = 1334
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 1214
$ws.Range("E163").Value = 38
$ws.Range("H163").Value = 82

# --- Row 167: Polinesia Francesa ----------------------------------------
$ws.Range("B167").Value = 1111
$ws.Range("D167").Value = 873
$ws.Range("E167").Value = 236

# --- Row 170: San Marino -------------------------------------------------
$ws.Range("D170").Value = 669
$ws.Range("E170").Value = 12

# --- Row 214: now Montserrat (was Islas Malvinas) -----------------------
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# --- Row 215: now Islas Malvinas (was Montserrat) -----------------------
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
